# Update the TestResults worksheet with the new iAuthor testcases.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

# Row 2: replace old "iAU_TC_ID_212" / "Pre-Request Verify Elumina..." test case
$ws.Range("A2").Value = "iAU_TC_ID_250"
$ws.Range("B2").Value = "@RegressionA Prerequisite Validation of Exam Approve`""
$ws.Range("C2").Value = "passed"

# Row 3: replace old "Pre-Request Validation of Delivery..." test case
$ws.Range("A3").Value = "iAU_TC_ID_250"
$ws.Range("B3").Value = "@RegressionA Validation of Exam Approve"
$ws.Range("C3").Value = "passed"

# Remove the two now-obsolete rows (previously rows 4 and 5)
$ws.Range("A4:C5").Delete()
